$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1939
$ws.Range("I9").Value = 2071.6667
$ws.Range("K9").Value = 2071.6667
$ws.Range("M9").Value = -1902.6667
$ws.Range("H12").Value = 391.69232
$ws.Range("I12").Value = 382.55554
$ws.Range("J12").Value = 412.25
$ws.Range("K12").Value = 382.55554
$ws.Range("L12").Value = 412.25
$ws.Range("M12").Value = -212.55554
$ws.Range("N12").Value = -752.25
$ws.Range("H17").Value = 4295265
$ws.Range("J17").Value = 4295265
$ws.Range("L17").Value = 12885795
$ws.Range("N17").Value = -12886131
$ws.Range("H28").Value = 768.3913
$ws.Range("I28").Value = 788.35
$ws.Range("K28").Value = 788.35
$ws.Range("M28").Value = -303.35
$ws.Range("H62").Value = 8357.200000000001
$ws.Range("I62").Value = 5122.2856
$ws.Range("K62").Value = 5122.2856
$ws.Range("M62").Value = -4498.2856
$ws.Range("H65").Value = 8357.200000000001
$ws.Range("I65").Value = 5122.2856
$ws.Range("K65").Value = 25611.428
$ws.Range("M65").Value = -22491.428
$ws.Range("H70").Value = 6988.6665
$ws.Range("J70").Value = 18499.666
$ws.Range("L70").Value = 55498.99800000001
$ws.Range("N70").Value = -56038.99800000001
$ws.Range("H73").Value = 6988.6665
$ws.Range("J73").Value = 18499.666
$ws.Range("L73").Value = 55498.99800000001
$ws.Range("N73").Value = -57370.99800000001
$ws.Range("H80").Value = 1050
$ws.Range("I80").Value = 1200
$ws.Range("J80").Value = 900
$ws.Range("K80").Value = 3600
$ws.Range("L80").Value = 2700
$ws.Range("M80").Value = -2602
$ws.Range("N80").Value = -4696
$ws.Range("H83").Value = 1050
$ws.Range("I83").Value = 1200
$ws.Range("J83").Value = 900
$ws.Range("K83").Value = 10800
$ws.Range("L83").Value = 8100
$ws.Range("M83").Value = -5808
$ws.Range("N83").Value = -18084
$ws.Range("H98").Value = 1704.0834
$ws.Range("I98").Value = 1720.9
$ws.Range("J98").Value = 1620
$ws.Range("K98").Value = 1720.9
$ws.Range("L98").Value = 1620
$ws.Range("M98").Value = -222.9000000000001
$ws.Range("N98").Value = -4616
$ws.Range("H113").Value = 17706.824
$ws.Range("I113").Value = 17101.215
$ws.Range("J113").Value = 20533
$ws.Range("K113").Value = 17101.215
$ws.Range("L113").Value = 20533
$ws.Range("M113").Value = -13847.215
$ws.Range("N113").Value = -27041
$ws.Range("H122").Value = 1704.0834
$ws.Range("I122").Value = 1720.9
$ws.Range("J122").Value = 1620
$ws.Range("K122").Value = 5162.700000000001
$ws.Range("L122").Value = 4860
$ws.Range("M122").Value = -2712.700000000001
$ws.Range("N122").Value = -9760
$ws.Range("H137").Value = 12138.794
$ws.Range("I137").Value = 3023.8635
$ws.Range("K137").Value = 9071.5905
$ws.Range("M137").Value = -6521.5905
$ws.Range("H138").Value = 2268.8918
$ws.Range("I138").Value = 2230.6924
$ws.Range("J138").Value = 2359.182
$ws.Range("K138").Value = 6692.0772
$ws.Range("L138").Value = 7077.545999999999
$ws.Range("M138").Value = -1552.0772
$ws.Range("N138").Value = -17357.546
$ws.Range("H141").Value = 4277
$ws.Range("I141").Value = 3880.5386
$ws.Range("J141").Value = 5995
$ws.Range("K141").Value = 11641.6158
$ws.Range("L141").Value = 17985
$ws.Range("M141").Value = -6461.6158
$ws.Range("N141").Value = -28345
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7584.21
$ws.Range("I32").Value = 2540.6125
$ws.Range("J32").Value = 27758.6
$ws.Range("K32").Value = 2540.6125
$ws.Range("L32").Value = 27758.6
$ws.Range("M32").Value = -2253.6125
$ws.Range("N32").Value = -28332.6
$ws.Range("H45").Value = 1578
$ws.Range("I45").Value = 1578
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1578
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1201
$ws.Range("N45").Value = $null
$ws.Range("H61").Value = 18742.363
$ws.Range("I61").Value = 19059.834
$ws.Range("J61").Value = 18623.312
$ws.Range("K61").Value = 19059.834
$ws.Range("L61").Value = 18623.312
$ws.Range("M61").Value = -18847.834
$ws.Range("N61").Value = -19047.312
$ws.Range("H122").Value = 2000
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = $null
$ws.Range("H136").Value = 18742.363
$ws.Range("I136").Value = 19059.834
$ws.Range("J136").Value = 18623.312
$ws.Range("K136").Value = 57179.50199999999
$ws.Range("L136").Value = 55869.936
$ws.Range("M136").Value = -54629.50199999999
$ws.Range("N136").Value = -60969.936
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 5115.8
$ws.Range("I94").Value = 4610.091
$ws.Range("K94").Value = 4610.091
$ws.Range("M94").Value = -4159.091
$ws.Range("H105").Value = 1364.8
$ws.Range("I105").Value = 1003.6
$ws.Range("K105").Value = 1003.6
$ws.Range("M105").Value = 743.4
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 280
$ws.Range("I2").Value = 152
$ws.Range("J2").Value = 600
$ws.Range("K2").Value = 152
$ws.Range("L2").Value = 600
$ws.Range("M2").Value = -39
$ws.Range("N2").Value = -826
$ws.Range("H16").Value = 9027.571
$ws.Range("J16").Value = 10383
$ws.Range("L16").Value = 10383
$ws.Range("N16").Value = -10957
$ws.Range("H58").Value = 13196.892
$ws.Range("I58").Value = 6164.154
$ws.Range("J58").Value = 17006.291
$ws.Range("K58").Value = 6164.154
$ws.Range("L58").Value = 17006.291
$ws.Range("M58").Value = -5961.154
$ws.Range("N58").Value = -17412.291
$ws.Range("H113").Value = 9027.571
$ws.Range("J113").Value = 10383
$ws.Range("L113").Value = 10383
$ws.Range("N113").Value = -14723
$ws.Range("H136").Value = 13196.892
$ws.Range("I136").Value = 6164.154
$ws.Range("J136").Value = 17006.291
$ws.Range("K136").Value = 18492.462
$ws.Range("L136").Value = 51018.87300000001
$ws.Range("M136").Value = -15942.462
$ws.Range("N136").Value = -56118.87300000001
$ws.Range("H141").Value = 226639.33
$ws.Range("J141").Value = 226639.33
$ws.Range("L141").Value = 226639.33
$ws.Range("N141").Value = -236999.33
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 4464798.5
$ws.Range("J107").Value = 10416927
$ws.Range("L107").Value = 31250781
$ws.Range("N107").Value = -31254621
$ws.Range("H108").Value = 483
$ws.Range("I108").Value = 483
$ws.Range("K108").Value = 1449
$ws.Range("M108").Value = 1431
$ws.Range("H109").Value = 5557004.5
$ws.Range("I109").Value = 1739
$ws.Range("J109").Value = 33333332
$ws.Range("K109").Value = 5217
$ws.Range("L109").Value = 99999996
$ws.Range("M109").Value = -4177
$ws.Range("N109").Value = -100002076
$ws.Range("H122").Value = 9997.5
$ws.Range("I122").Value = 1795.4
$ws.Range("J122").Value = 15856.143
$ws.Range("K122").Value = 16158.6
$ws.Range("L122").Value = 142705.287
$ws.Range("M122").Value = -13708.6
$ws.Range("N122").Value = -147605.287
$ws.Range("H131").Value = 1483.52
$ws.Range("J131").Value = 1493.3877
$ws.Range("L131").Value = 4480.1631
$ws.Range("N131").Value = -14560.1631
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1200.2
$ws.Range("I107").Value = 849.5
$ws.Range("J107").Value = 1434
$ws.Range("K107").Value = 849.5
$ws.Range("L107").Value = 1434
$ws.Range("M107").Value = 1070.5
$ws.Range("N107").Value = -5274
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = $null
$ws.Range("N122").Value = $null
$ws.Range("H126").Value = 7553.3667
$ws.Range("J126").Value = 7755.1665
$ws.Range("L126").Value = 23265.4995
$ws.Range("N126").Value = -28205.4995
$ws.Range("H132").Value = 6248.5674
$ws.Range("I132").Value = 6248.5674
$ws.Range("K132").Value = 18745.7022
$ws.Range("M132").Value = -16215.7022
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 21252688
$ws.Range("I2").Value = 3666.6667
$ws.Range("J2").Value = 34002100
$ws.Range("K2").Value = 3666.6667
$ws.Range("L2").Value = 34002100
$ws.Range("M2").Value = -3554.6667
$ws.Range("N2").Value = -34002324
$ws.Range("H40").Value = 10024.417
$ws.Range("I40").Value = 1616.5
$ws.Range("K40").Value = 1616.5
$ws.Range("M40").Value = -1480.5
$ws.Range("H61").Value = 4474.15
$ws.Range("I61").Value = 2903.3
$ws.Range("K61").Value = 2903.3
$ws.Range("M61").Value = -2701.3
$ws.Range("H113").Value = 4474.15
$ws.Range("I113").Value = 2903.3
$ws.Range("K113").Value = 2903.3
$ws.Range("M113").Value = -733.3000000000002
$ws.Range("H136").Value = 13993
$ws.Range("I136").Value = 15682.667
$ws.Range("K136").Value = 47048.001
$ws.Range("M136").Value = -44498.001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 9000
$ws.Range("I29").Value = 9250
$ws.Range("J29").Value = 8500
$ws.Range("K29").Value = 9250
$ws.Range("L29").Value = 8500
$ws.Range("M29").Value = -8960
$ws.Range("N29").Value = -9080
$ws.Range("H64").Value = 52962
$ws.Range("J64").Value = 52966.43
$ws.Range("L64").Value = 52966.43
$ws.Range("N64").Value = -53462.43
$ws.Range("H67").Value = 52962
$ws.Range("J67").Value = 52966.43
$ws.Range("L67").Value = 52966.43
$ws.Range("N67").Value = -54682.43
$ws.Range("H122").Value = 4692.2974
$ws.Range("I122").Value = 2361.2173
$ws.Range("K122").Value = 7083.651899999999
$ws.Range("M122").Value = -4633.651899999999
$ws.Range("H136").Value = 14748.895
$ws.Range("I136").Value = 2886.4
$ws.Range("K136").Value = 8659.200000000001
$ws.Range("M136").Value = -6109.200000000001
